# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# Update the Metadata sheet: version bump, new publish date, publisher, and
# jurisdiction; also drop the stray duplicate "Contact" row.
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Cells.Item(3, 2).Value = "6.0.0"
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"
$meta.Cells.Item(9, 2).Value = "Alvearie Team"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# remove it entirely so everything below shifts up by one row.
$meta.Rows.Item(11).Delete()

# Update the Elements sheet: the root "Extension" element's Short/Definition
# text was placeholder boilerplate; replace it with the real category text.
$elements = $wb.Worksheets.Item("Elements")
$elements.Cells.Item(2, 11).Value = "Category"
$elements.Cells.Item(2, 12).Value = "Top level category for classification purposes"
